# Insert a new data row at row 154, shifting existing rows 154..262 down to 155..263.
# This mirrors the authored change: a new price quote record was added to the
# "Membrillo" (Vega Modelo de Temuco) dataset, pushing every subsequent row down
# by one position (dimension grows from A1:T262 to A1:T263).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Insert()

$ws.Cells.Item(154, 1).Value = 10
$ws.Cells.Item(154, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(154, 3).Value = "La Araucanía"
$ws.Cells.Item(154, 4).Value = 45040
$ws.Cells.Item(154, 5).Value = 9
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100104
$ws.Cells.Item(154, 8).Value = "Frutos de pepita"
$ws.Cells.Item(154, 9).Value = 100104003
$ws.Cells.Item(154, 10).Value = "Membrillo"
$ws.Cells.Item(154, 11).Value = "Champion"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 110
$ws.Cells.Item(154, 14).Value = 13000
$ws.Cells.Item(154, 15).Value = 13000
$ws.Cells.Item(154, 16).Value = 13000
$ws.Cells.Item(154, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(154, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(154, 19).Value = 722
$ws.Cells.Item(154, 20).Value = 18
